$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 25725.818
$ws.Range("I64").Value = 29609.389
$ws.Range("J64").Value = 8249.75
$ws.Range("K64").Value = 29609.389
$ws.Range("L64").Value = 8249.75
$ws.Range("M64").Value = -29361.389
$ws.Range("N64").Value = -8745.75
$ws.Range("H67").Value = 25725.818
$ws.Range("I67").Value = 29609.389
$ws.Range("J67").Value = 8249.75
$ws.Range("K67").Value = 29609.389
$ws.Range("L67").Value = 8249.75
$ws.Range("M67").Value = -28751.389
$ws.Range("N67").Value = -9965.75
$ws.Range("H110").Value = 702000000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 702000000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 702000000
$ws.Range("N110").Value = -702008180
$ws.Range("H111").Value = 2157.3333
$ws.Range("I111").Value = 2063.8
$ws.Range("J111").Value = 2274.25
$ws.Range("K111").Value = 6191.400000000001
$ws.Range("L111").Value = 6822.75
$ws.Range("M111").Value = -3124.400000000001
$ws.Range("N111").Value = -12956.75
$ws.Range("H113").Value = 11741.385
$ws.Range("J113").Value = 6625.4
$ws.Range("L113").Value = 6625.4
$ws.Range("N113").Value = -13133.4
$ws.Range("H137").Value = 16760.268
$ws.Range("I137").Value = 21137.727
$ws.Range("J137").Value = 4722.25
$ws.Range("K137").Value = 63413.181
$ws.Range("L137").Value = 14166.75
$ws.Range("M137").Value = -60863.181
$ws.Range("N137").Value = -19266.75
$ws.Range("H138").Value = 3416.261
$ws.Range("I138").Value = 546.3570999999999
$ws.Range("J138").Value = 4146.7817
$ws.Range("K138").Value = 1639.0713
$ws.Range("L138").Value = 12440.3451
$ws.Range("M138").Value = 3500.9287
$ws.Range("N138").Value = -22720.3451
$ws.Range("M110").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 12250
$ws.Range("I39").Value = 6750
$ws.Range("K39").Value = 6750
$ws.Range("M39").Value = -6230
$ws.Range("H45").Value = 138663.27
$ws.Range("I45").Value = 204081.9
$ws.Range("K45").Value = 204081.9
$ws.Range("M45").Value = -203704.9
$ws.Range("H61").Value = 7338.3076
$ws.Range("I61").Value = 8376.777
$ws.Range("K61").Value = 8376.777
$ws.Range("M61").Value = -8164.777
$ws.Range("H122").Value = 414275.75
$ws.Range("I122").Value = 1851
$ws.Range("J122").Value = 2338924.8
$ws.Range("K122").Value = 5553
$ws.Range("L122").Value = 7016774.399999999
$ws.Range("M122").Value = -3103
$ws.Range("N122").Value = -7021674.399999999
$ws.Range("H132").Value = 1724
$ws.Range("I132").Value = 915.25714
$ws.Range("K132").Value = 2745.77142
$ws.Range("M132").Value = -215.77142
$ws.Range("H136").Value = 7338.3076
$ws.Range("I136").Value = 8376.777
$ws.Range("K136").Value = 25130.331
$ws.Range("M136").Value = -22580.331

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 5629.722
$ws.Range("I64").Value = 8624.777
$ws.Range("J64").Value = 2634.6667
$ws.Range("K64").Value = 8624.777
$ws.Range("L64").Value = 2634.6667
$ws.Range("M64").Value = -8399.777
$ws.Range("N64").Value = -3084.6667
$ws.Range("H67").Value = 5629.722
$ws.Range("I67").Value = 8624.777
$ws.Range("J67").Value = 2634.6667
$ws.Range("K67").Value = 8624.777
$ws.Range("L67").Value = 2634.6667
$ws.Range("M67").Value = -7844.777
$ws.Range("N67").Value = -4194.6667
$ws.Range("H134").Value = 2546.925
$ws.Range("I134").Value = 1652.3529
$ws.Range("K134").Value = 4957.0587
$ws.Range("M134").Value = -2422.0587

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 551.2593000000001
$ws.Range("I22").Value = 515.1429000000001
$ws.Range("J22").Value = 563.9
$ws.Range("K22").Value = 515.1429000000001
$ws.Range("L22").Value = 563.9
$ws.Range("M22").Value = -165.1429000000001
$ws.Range("N22").Value = -1263.9
$ws.Range("H31").Value = 4818
$ws.Range("I31").Value = 3384.4285
$ws.Range("J31").Value = 5933
$ws.Range("K31").Value = 3384.4285
$ws.Range("L31").Value = 5933
$ws.Range("M31").Value = -3089.4285
$ws.Range("N31").Value = -6523
$ws.Range("H34").Value = 4818
$ws.Range("I34").Value = 3384.4285
$ws.Range("J34").Value = 5933
$ws.Range("K34").Value = 3384.4285
$ws.Range("L34").Value = 5933
$ws.Range("M34").Value = -3182.4285
$ws.Range("N34").Value = -6337
$ws.Range("H86").Value = 10970.786
$ws.Range("I86").Value = 8732.333000000001
$ws.Range("K86").Value = 8732.333000000001
$ws.Range("M86").Value = -7609.333000000001
$ws.Range("H89").Value = 10970.786
$ws.Range("I89").Value = 8732.333000000001
$ws.Range("K89").Value = 43661.665
$ws.Range("M89").Value = -38045.665
$ws.Range("H99").Value = 11626046
$ws.Range("J99").Value = 3850
$ws.Range("L99").Value = 3850
$ws.Range("N99").Value = -6846
$ws.Range("H107").Value = 16664
$ws.Range("I107").Value = 24874.445
$ws.Range("K107").Value = 24874.445
$ws.Range("M107").Value = -22954.445
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H126").Value = 11626046
$ws.Range("J126").Value = 3850
$ws.Range("L126").Value = 11550
$ws.Range("N126").Value = -16490
$ws.Range("N111").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1213.75
$ws.Range("I38").Value = 323.91666
$ws.Range("K38").Value = 971.7499799999999
$ws.Range("M38").Value = -624.7499799999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5955.5
$ws.Range("I102").Value = 7045.48
$ws.Range("K102").Value = 7045.48
$ws.Range("M102").Value = -5423.48
$ws.Range("H122").Value = 17676.182
$ws.Range("I122").Value = 14099.286
$ws.Range("K122").Value = 42297.858
$ws.Range("M122").Value = -39847.858
$ws.Range("H126").Value = 11990.728
$ws.Range("I126").Value = 9329.143
$ws.Range("J126").Value = 16648.5
$ws.Range("K126").Value = 27987.429
$ws.Range("L126").Value = 49945.5
$ws.Range("M126").Value = -25517.429
$ws.Range("N126").Value = -54885.5
$ws.Range("H132").Value = 4287.4443
$ws.Range("I132").Value = 4056.4285
$ws.Range("J132").Value = 5096
$ws.Range("K132").Value = 12169.2855
$ws.Range("L132").Value = 15288
$ws.Range("M132").Value = -9639.2855
$ws.Range("N132").Value = -20348

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 203999
$ws.Range("I7").Value = 203999
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 203999
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -203887
$ws.Range("H40").Value = 44239.25
$ws.Range("I40").Value = 53498.375
$ws.Range("J40").Value = 25721
$ws.Range("K40").Value = 53498.375
$ws.Range("L40").Value = 25721
$ws.Range("M40").Value = -53362.375
$ws.Range("N40").Value = -25993
$ws.Range("H122").Value = 3314.8
$ws.Range("I122").Value = 3314.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9944.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7494.400000000001
$ws.Range("H126").Value = 203999
$ws.Range("I126").Value = 203999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 611997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -609527
$ws.Range("H132").Value = 881732.0600000001
$ws.Range("I132").Value = 2484909.2
$ws.Range("K132").Value = 7454727.600000001
$ws.Range("M132").Value = -7452197.600000001
$ws.Range("H136").Value = 7743.7
$ws.Range("I136").Value = 2713.8572
$ws.Range("J136").Value = 10452.077
$ws.Range("K136").Value = 8141.571599999999
$ws.Range("L136").Value = 31356.231
$ws.Range("M136").Value = -5591.571599999999
$ws.Range("N136").Value = -36456.231
$ws.Range("N7").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("N126").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2330
$ws.Range("I14").Value = 2330
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2330
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2162
$ws.Range("H107").Value = 45856.855
$ws.Range("I107").Value = 3499.6667
$ws.Range("J107").Value = 300000
$ws.Range("K107").Value = 10499.0001
$ws.Range("L107").Value = 900000
$ws.Range("M107").Value = -8579.000100000001
$ws.Range("N107").Value = -903840
$ws.Range("H122").Value = 18152.305
$ws.Range("I122").Value = 3362.3704
$ws.Range("J122").Value = 62522.11
$ws.Range("K122").Value = 10087.1112
$ws.Range("L122").Value = 187566.33
$ws.Range("M122").Value = -7637.111199999999
$ws.Range("N122").Value = -192466.33
$ws.Range("N14").ClearContents()
